$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H38").Value = 712
$ws.Range("I38").Value = 712
$ws.Range("K38").Value = 2136
$ws.Range("M38").Value = -1764

$ws.Range("H58").Value = 500
$ws.Range("J58").Value = 900
$ws.Range("L58").Value = 2700
$ws.Range("N58").Value = -3000

$ws.Range("H61").Value = 351.25
$ws.Range("I61").Value = 351.25
$ws.Range("K61").Value = 1053.75
$ws.Range("M61").Value = -881.75

$ws.Range("H100").Value = 5462.8335
$ws.Range("J100").Value = 6647.5835
$ws.Range("L100").Value = 6647.5835
$ws.Range("N100").Value = -7729.5835

$ws.Range("H137").Value = 2343.15
$ws.Range("I137").Value = 1641.0667
$ws.Range("J137").Value = 4449.4
$ws.Range("K137").Value = 4923.2001
$ws.Range("L137").Value = 13348.2
$ws.Range("M137").Value = -2373.2001
$ws.Range("N137").Value = -18448.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11878.104
$ws.Range("I2").Value = 12825.577
$ws.Range("J2").Value = 3666.6667
$ws.Range("K2").Value = 12825.577
$ws.Range("L2").Value = 3666.6667
$ws.Range("M2").Value = -12712.577
$ws.Range("N2").Value = -3892.6667

$ws.Range("H45").Value = 6166.636
$ws.Range("I45").Value = 4763.8
$ws.Range("K45").Value = 4763.8
$ws.Range("M45").Value = -4386.8

$ws.Range("H61").Value = 4507.4585
$ws.Range("I61").Value = 3230.0908
$ws.Range("K61").Value = 3230.0908
$ws.Range("M61").Value = -3018.0908

$ws.Range("H74").Value = 2232.423
$ws.Range("I74").Value = 1957.174
$ws.Range("K74").Value = 1957.174
$ws.Range("M74").Value = -1083.174

$ws.Range("H77").Value = 2232.423
$ws.Range("I77").Value = 1957.174
$ws.Range("K77").Value = 9785.869999999999
$ws.Range("M77").Value = -5417.869999999999

$ws.Range("H102").Value = 6000
$ws.Range("I102").Value = 6000
$ws.Range("K102").Value = 6000
$ws.Range("M102").Value = -4378

$ws.Range("H116").Value = 11878.104
$ws.Range("I116").Value = 12825.577
$ws.Range("J116").Value = 3666.6667
$ws.Range("K116").Value = 12825.577
$ws.Range("L116").Value = 3666.6667
$ws.Range("M116").Value = -10531.577
$ws.Range("N116").Value = -8254.6667

$ws.Range("H122").Value = 5041.5815
$ws.Range("I122").Value = 4295.724
$ws.Range("J122").Value = 6586.5713
$ws.Range("K122").Value = 12887.172
$ws.Range("L122").Value = 19759.7139
$ws.Range("M122").Value = -10437.172
$ws.Range("N122").Value = -24659.7139

$ws.Range("H132").Value = 34487156
$ws.Range("I132").Value = 50003750
$ws.Range("K132").Value = 150011250
$ws.Range("M132").Value = -150008720

$ws.Range("H136").Value = 4507.4585
$ws.Range("I136").Value = 3230.0908
$ws.Range("K136").Value = 9690.2724
$ws.Range("M136").Value = -7140.2724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11878.104
$ws.Range("I3").Value = 12825.577
$ws.Range("J3").Value = 3666.6667
$ws.Range("K3").Value = 12825.577
$ws.Range("L3").Value = 3666.6667
$ws.Range("M3").Value = -12711.577
$ws.Range("N3").Value = -3894.6667

$ws.Range("H86").Value = 5973.909
$ws.Range("I86").Value = 1583.3334
$ws.Range("J86").Value = 7620.375
$ws.Range("K86").Value = 1583.3334
$ws.Range("L86").Value = 7620.375
$ws.Range("M86").Value = -460.3334
$ws.Range("N86").Value = -9866.375

$ws.Range("H89").Value = 5973.909
$ws.Range("I89").Value = 1583.3334
$ws.Range("J89").Value = 7620.375
$ws.Range("K89").Value = 7916.666999999999
$ws.Range("L89").Value = 38101.875
$ws.Range("M89").Value = -2300.666999999999
$ws.Range("N89").Value = -49333.875

$ws.Range("H128").Value = 3958.1667
$ws.Range("I128").Value = 3958.1667
$ws.Range("K128").Value = 11874.5001
$ws.Range("M128").Value = -9384.500100000001

$ws.Range("H134").Value = 4214.5713
$ws.Range("I134").Value = 2683
$ws.Range("J134").Value = 6971.4
$ws.Range("K134").Value = 8049
$ws.Range("L134").Value = 20914.2
$ws.Range("M134").Value = -5514
$ws.Range("N134").Value = -25984.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3958.7942
$ws.Range("I31").Value = 965.2857
$ws.Range("J31").Value = 6054.25
$ws.Range("K31").Value = 965.2857
$ws.Range("L31").Value = 6054.25
$ws.Range("M31").Value = -670.2857
$ws.Range("N31").Value = -6644.25

$ws.Range("H34").Value = 3958.7942
$ws.Range("I34").Value = 965.2857
$ws.Range("J34").Value = 6054.25
$ws.Range("K34").Value = 965.2857
$ws.Range("L34").Value = 6054.25
$ws.Range("M34").Value = -763.2857
$ws.Range("N34").Value = -6458.25

$ws.Range("H58").Value = 4159.5
$ws.Range("I58").Value = 2649.5
$ws.Range("J58").Value = 4662.8335
$ws.Range("K58").Value = 2649.5
$ws.Range("L58").Value = 4662.8335
$ws.Range("M58").Value = -2446.5
$ws.Range("N58").Value = -5068.8335

$ws.Range("H99").Value = 21882594
$ws.Range("I99").Value = 25005124
$ws.Range("K99").Value = 25005124
$ws.Range("M99").Value = -25003626

$ws.Range("H126").Value = 21882594
$ws.Range("I126").Value = 25005124
$ws.Range("K126").Value = 75015372
$ws.Range("M126").Value = -75012902

$ws.Range("H134").Value = 4819.119
$ws.Range("I134").Value = 3636.6296
$ws.Range("J134").Value = 6947.6
$ws.Range("K134").Value = 10909.8888
$ws.Range("L134").Value = 20842.8
$ws.Range("M134").Value = -8374.888800000001
$ws.Range("N134").Value = -25912.8

$ws.Range("H136").Value = 4159.5
$ws.Range("I136").Value = 2649.5
$ws.Range("J136").Value = 4662.8335
$ws.Range("K136").Value = 7948.5
$ws.Range("L136").Value = 13988.5005
$ws.Range("M136").Value = -5398.5
$ws.Range("N136").Value = -19088.5005

$ws.Range("H141").Value = 91903.19500000001
$ws.Range("J141").Value = 94733.3
$ws.Range("L141").Value = 94733.3
$ws.Range("N141").Value = -105093.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 139026.53
$ws.Range("I63").Value = 501493.75
$ws.Range("J63").Value = 7220.273
$ws.Range("K63").Value = 1504481.25
$ws.Range("L63").Value = 21660.819
$ws.Range("M63").Value = -1503732.25
$ws.Range("N63").Value = -23158.819

$ws.Range("H64").Value = 3180.9092
$ws.Range("I64").Value = 4998
$ws.Range("K64").Value = 14994
$ws.Range("M64").Value = -14724

$ws.Range("H66").Value = 139026.53
$ws.Range("I66").Value = 501493.75
$ws.Range("J66").Value = 7220.273
$ws.Range("K66").Value = 4513443.75
$ws.Range("L66").Value = 64982.457
$ws.Range("M66").Value = -4509699.75
$ws.Range("N66").Value = -72470.45699999999

$ws.Range("H67").Value = 3180.9092
$ws.Range("I67").Value = 4998
$ws.Range("K67").Value = 14994
$ws.Range("M67").Value = -14058

$ws.Range("H97").Value = 967.96
$ws.Range("J97").Value = 473.46667
$ws.Range("L97").Value = 1420.40001
$ws.Range("N97").Value = -2412.40001

$ws.Range("H117").Value = 565
$ws.Range("I117").Value = 128
$ws.Range("J117").Value = 843.0909
$ws.Range("K117").Value = 384
$ws.Range("L117").Value = 2529.2727
$ws.Range("M117").Value = 3058
$ws.Range("N117").Value = -9413.2727

$ws.Range("H122").Value = 6290364.5
$ws.Range("I122").Value = 851.4138
$ws.Range("K122").Value = 7662.724200000001
$ws.Range("M122").Value = -5212.724200000001

$ws.Range("H131").Value = 5845.9375
$ws.Range("I131").Value = 3717.1428
$ws.Range("J131").Value = 7501.6665
$ws.Range("K131").Value = 11151.4284
$ws.Range("L131").Value = 22504.9995
$ws.Range("M131").Value = -6111.428400000001
$ws.Range("N131").Value = -32584.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3910.2812
$ws.Range("I132").Value = 2890.318
$ws.Range("K132").Value = 8670.954000000002
$ws.Range("M132").Value = -6140.954000000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 11224.158
$ws.Range("I40").Value = 15775.125
$ws.Range("K40").Value = 15775.125
$ws.Range("M40").Value = -15639.125

$ws.Range("H136").Value = 4787.5674
$ws.Range("I136").Value = 3549.5454
$ws.Range("K136").Value = 10648.6362
$ws.Range("M136").Value = -8098.636200000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 999999
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H62").Value = 3775
$ws.Range("I62").Value = 100
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 100
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = 524
$ws.Range("N62").Value = -6248

$ws.Range("H65").Value = 3775
$ws.Range("I65").Value = 100
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 500
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = 2620
$ws.Range("N65").Value = -31240

$ws.Range("H81").Value = 10104.667
$ws.Range("J81").Value = 4750
$ws.Range("L81").Value = 9500
$ws.Range("N81").Value = -11622

$ws.Range("H84").Value = 10104.667
$ws.Range("J84").Value = 4750
$ws.Range("L84").Value = 47500
$ws.Range("N84").Value = -58108

$ws.Range("H136").Value = 4130.5
$ws.Range("I136").Value = 3480.6086
$ws.Range("K136").Value = 10441.8258
$ws.Range("M136").Value = -7891.825800000001
